$d = $word.ActiveDocument

# --- Locate the target paragraph ---
# "Main Page" occurs as its own run in three places in the document, but only
# the second occurrence is the paragraph we need (the one split into
# "...Vehicles" page which is in the " / "Main Page" / " of the of the website." runs).
$firstMainPage = $d.Content
$firstMainPage.Find.Execute("Main Page", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$afterFirst = $d.Range($firstMainPage.End, $d.Content.End)
$afterFirst.Find.Execute("Main Page", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$mainPageRange = $d.Range($afterFirst.Start, $afterFirst.End)

# --- Fix 1: " of the of the website." -> " of the website." (only this paragraph) ---
$tail = $d.Range($mainPageRange.End, $mainPageRange.End + 40)
$tail.Find.Execute(" of the of the website.", $true, $false, $false, $false, $false, $true, 1, $false, " of the website.", 2)

# --- Fix 2: split the run right before "page which is in the" and move the
# _GoBack bookmark there (Bookmarks.Add relocates it, removing the old one). ---
$before = $d.Range($mainPageRange.Start - 200, $mainPageRange.Start)
$before.Find.Execute("page which is in the", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$splitPoint = $d.Range($before.Start, $before.Start)
$d.Bookmarks.Add("_GoBack", $splitPoint)
